# Sync attendance_reports, modules_schedules, and assets from main repo - 2025-10-27 06:24:37
#
# Applies the session-analysis corrections to the "Session Analysis Results" sheet:
#  - widen the Status column (I) slightly
#  - re-sync the "Recorded By" grader lists for several sessions
#  - refresh the class/group attendance statistics
#  - record the HISTOLOGY A3 session #1 as "Not Recorded" (new pink status style)
#  - record the ANATOMY B1 session #2 as newly "Recorded"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Column I (Status) width 10 -> 14
# ---------------------------------------------------------------------------
$ws.Columns.Item(9).ColumnWidth = 13.2

# ---------------------------------------------------------------------------
# 2. Re-synced "Recorded By" grader lists (order changed upstream; same people)
# ---------------------------------------------------------------------------
$gradersA = "nesmadrahim@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, servinaz@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg"
$ws.Range("G2").Value = $gradersA
$ws.Range("G17").Value = $gradersA
$ws.Range("G92").Value = $gradersA
$ws.Range("G107").Value = $gradersA

$gradersB = "nahla.nagiub@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, servinaz@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, gehanadel@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"
$ws.Range("G32").Value = $gradersB
$ws.Range("G47").Value = $gradersB

$gradersC = "rana.abozaid@med.asu.edu.eg, servinaz@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg"
$ws.Range("G33").Value = $gradersC
$ws.Range("G48").Value = $gradersC

$gradersD = "hend_mahmoud@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, servinaz@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, gehanadel@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg"
$ws.Range("G62").Value = $gradersD
$ws.Range("G77").Value = $gradersD

# ---------------------------------------------------------------------------
# 3. Class Statistics block (rows 6-10)
# ---------------------------------------------------------------------------
$ws.Range("L6").Value = 17
$ws.Range("L7").Value = 1
$ws.Range("L8").Value = 102
$ws.Range("L9").Value = "14.2%"
$ws.Range("L10").Value = "49.1%"

# ---------------------------------------------------------------------------
# 4. Group Statistics block (rows 17, 19)
# ---------------------------------------------------------------------------
$ws.Range("P17").Value = 1
$ws.Range("Q17").Value = 12

$ws.Range("O19").Value = 3
$ws.Range("Q19").Value = 12
$ws.Range("R19").Value = "20.0%"
$ws.Range("S19").Value = "42.9%"

# ---------------------------------------------------------------------------
# 5. Row 37 (Year 2 / A3 / HISTOLOGY #1): now "Not Recorded" -> new pink style
#    Build the style by cloning an existing centered cell (keeps the shared
#    font/alignment) and then recoloring its fill to match the "Not Recorded"
#    legend color (light pink, FFB6C1) used elsewhere in the workbook.
# ---------------------------------------------------------------------------
$row37 = $ws.Range("A37:I37")
$ws.Range("A2").Copy($row37)
$row37.Interior.Pattern = 1
$row37.Interior.Color = 12695295
$ws.Range("I37").Value = "Not Recorded"

# ---------------------------------------------------------------------------
# 6. Row 63 (Year 2 / B1 / ANATOMY #2): now "Recorded" -> reuse the standard
#    "Recorded" style (same as row 62 / style used by A2) and fill in results.
# ---------------------------------------------------------------------------
$row63 = $ws.Range("A63:I63")
$ws.Range("A2").Copy($row63)
$ws.Range("G63").Value = "mennatulla.medhat@med.asu.edu.eg"
$ws.Range("H63").Value = "81/154"
$ws.Range("I63").Value = "Recorded"
